{"js": "// Auto-generated replacements for three-digit/one-digit division problems.\n// Each entry is [oldText, newText]; oldText is unique across the document body,\n// so a direct search+replace (matchCase, exact) is safe and order-independent.\nconst replacements = [\n  [\"947\u00f79=\", \"325\u00f76=\"],\n  [\"245\u00f77=\", \"740\u00f78=\"],\n  [\"847\u00f73=\", \"498\u00f76=\"],\n  [\"754\u00f76=\", \"889\u00f72=\"],\n  [\"489\u00f74=\", \"555\u00f72=\"],\n  [\"969\u00f77=\", \"297\u00f78=\"],\n  [\"897\u00f76=\", \"728\u00f76=\"],\n  [\"685\u00f76=\", \"672\u00f75=\"],\n  [\"747\u00f76=\", \"278\u00f75=\"],\n  [\"816\u00f73=\", \"867\u00f76=\"],\n  [\"983\u00f75=\", \"432\u00f78=\"],\n  [\"680\u00f78=\", \"800\u00f72=\"],\n  [\"637\u00f72=\", \"482\u00f76=\"],\n  [\"914\u00f78=\", \"557\u00f74=\"],\n  [\"370\u00f78=\", \"262\u00f77=\"],\n  [\"739\u00f72=\", \"622\u00f77=\"],\n  [\"852\u00f79=\", \"457\u00f74=\"],\n  [\"852\u00f77=\", \"849\u00f75=\"],\n  [\"246\u00f72=\", \"319\u00f78=\"],\n  [\"693\u00f74=\", \"958\u00f77=\"],\n  [\"431\u00f77=\", \"526\u00f75=\"],\n  [\"973\u00f79=\", \"252\u00f77=\"],\n  [\"929\u00f76=\", \"683\u00f79=\"],\n  [\"425\u00f79=\", \"598\u00f76=\"],\n  [\"975\u00f72=\", \"421\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, 'Replace');\n  }\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated replacements for three-digit/one-digit division problems.\n# Each entry is old/new text; old text is unique across the document body,\n# so Find/Replace against the whole-document Range is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('947\u00f79=', '325\u00f76='),\n    @('245\u00f77=', '740\u00f78='),\n    @('847\u00f73=', '498\u00f76='),\n    @('754\u00f76=', '889\u00f72='),\n    @('489\u00f74=', '555\u00f72='),\n    @('969\u00f77=', '297\u00f78='),\n    @('897\u00f76=', '728\u00f76='),\n    @('685\u00f76=', '672\u00f75='),\n    @('747\u00f76=', '278\u00f75='),\n    @('816\u00f73=', '867\u00f76='),\n    @('983\u00f75=', '432\u00f78='),\n    @('680\u00f78=', '800\u00f72='),\n    @('637\u00f72=', '482\u00f76='),\n    @('914\u00f78=', '557\u00f74='),\n    @('370\u00f78=', '262\u00f77='),\n    @('739\u00f72=', '622\u00f77='),\n    @('852\u00f79=', '457\u00f74='),\n    @('852\u00f77=', '849\u00f75='),\n    @('246\u00f72=', '319\u00f78='),\n    @('693\u00f74=', '958\u00f77='),\n    @('431\u00f77=', '526\u00f75='),\n    @('973\u00f79=', '252\u00f77='),\n    @('929\u00f76=', '683\u00f79='),\n    @('425\u00f79=', '598\u00f76='),\n    @('975\u00f72=', '421\u00f78='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Search text not found: $oldText\"\n    }\n}\n"}
